# Applies:
#  1) Bump every cached "datetimeFigureOut" auto-date placeholder
#     (slide master, every slide layout, handout master, notes master)
#     from 12/28/2022 to 12/29/2022.
#  2) Re-word the title & body copy on slide 1 ("Essential Presentation"
#     -> "Adventure Works Cycles" boilerplate).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "12/29/2022"

# Slide master
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Handout master
Set-DatePlaceholderText $p.HandoutMaster.Shapes $newDate

# Notes master
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# ---------------------------------------------------------------------
# Slide 1 content edits
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# Title placeholder: "Essential Presentation" -> "Adventure Works Cycles"
$titleShape = $slide1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Replace("Essential Presentation", "Adventure Works Cycles") | Out-Null

# Subtitle/body placeholder: replace both paragraphs with the new copy
$bodyShape = $slide1.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

$para1 = "Adventure Works Cycles, the fictitious company on which the Adventure Works sample databases are based, is a large, multinational manufacturing company. The company manufactures and sells metal and composite bicycles to North American, European and Asian commercial markets. "
$para2 = "In 2000, Adventure Works Cycles bought a small manufacturing plant, Importadores Neptuno, located in Mexico. Importadores Neptuno manufactures several critical subcomponents for the Adventure Works Cycles product line. These subcomponents are shipped to the Bothell location for final product assembly. "

$bodyRange.Paragraphs(2, 1).Text = $para1
$bodyRange.Paragraphs(3, 1).Text = $para2
